# Append new listings and refresh existing ones on the "ランサーズ" sheet,
# matching the commit "Append: 2026-01-14 18:30 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-01-14 18:30:57"

$rows = @(
  @{ Row = 2; Title = "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)"; Category = "システム開発"; Price = "200,000 円 ~ 300,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5450864"; Score = 383; Skill = "🔥AI,Ai ◆開発" },
  @{ Row = 3; Title = "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集"; Category = "システム開発"; Price = "1,000,000 円 ~ 3,000,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5460294"; Score = 375; Skill = "🔥AI,Ai ◆開発" },
  @{ Row = 4; Title = "【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集"; Category = "システム開発"; Price = "500,000 円 ~ 1,000,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5460267"; Score = 375; Skill = "🔥AI,Ai ◆開発" },
  @{ Row = 5; Title = "初回 資格学習支援ボット開発(Make/Airtable/Gemini/Stripe連携)LINE×AI"; Category = "システム開発"; Price = "200,000 円 ~ 300,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5471108"; Score = 368; Skill = "🔥AI,Ai ◆開発" },
  @{ Row = 6; Title = "AI企画書作成システムの「見積書作成」をご支援いただける制作会社/エンジニア募集(発注確約なし)"; Category = "システム開発"; Price = "20,000 円 ~ 50,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470737"; Score = 313; Skill = "🔥AI,Ai" },
  @{ Row = 7; Title = "AIオペレーションデザイナーを募集します(経営直結/会議→意思決定変換)"; Category = "システム開発"; Price = "300,000 円 ~ 500,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5471032"; Score = 310; Skill = "🔥AI,Ai" },
  @{ Row = 8; Title = "googleビジネスプロフィール一括ツール"; Category = "システム開発"; Price = "50,000 円 ~ 100,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470814"; Score = 73; Skill = "◆ツール" },
  @{ Row = 9; Title = "エクセルを利用して売上管理"; Category = "システム開発"; Price = "50,000 円 ~ 100,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5471068"; Score = 38; Skill = "◇管理" },
  @{ Row = 10; Title = "進行管理およびチームディレクションを担当"; Category = "システム開発"; Price = "~ 5,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5418064"; Score = 30; Skill = "◇管理" },
  @{ Row = 11; Title = "Rubyの暗号化機能のPHP化"; Category = "システム開発"; Price = "20,000 円 ~ 50,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470623"; Score = 28; Skill = "○PHP" },
  @{ Row = 12; Title = "オンラインWEB予約システム構築"; Category = "システム開発"; Price = "100,000 円 ~ 200,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470812"; Score = 33; Skill = "" },
  @{ Row = 13; Title = "金融機関の入出金伝票印刷システム構築依頼"; Category = "システム開発"; Price = "20,000 円 ~ 50,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470403"; Score = 28; Skill = "" },
  @{ Row = 14; Title = "移動型筐体の制御ハーネス製作(Arduino/電飾/音声/電源)"; Category = "システム開発"; Price = "50,000 円 ~ 100,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5471022"; Score = 18; Skill = "" },
  @{ Row = 15; Title = "《長期レギュラー》公的機関Web運用の要となる、ディレクター募集"; Category = "システム開発"; Price = "200,000 円 ~ 300,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470150"; Score = 18; Skill = "" },
  @{ Row = 16; Title = "限定公開 限定公開の仕事"; Category = "システム開発"; Price = "20,000 円 ~ 50,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5450323"; Score = 13; Skill = "" },
  @{ Row = 17; Title = "【フリーランス必見】エンジニア支援サービスのご紹介!"; Category = "システム開発"; Price = "10,000 円 ~ 20,000 円 / 固定"; Deadline = "期限情報なし"; Url = "https://www.lancers.jp/work/detail/5470726"; Score = 10; Skill = "" }
)

# Remove all existing hyperlinks first (they will be rebuilt below, since the
# row each one refers to is shifting as new listings are inserted at the top).
$ws.Hyperlinks.Delete()

# Clear out the old data rows completely; the freshest dataset (16 listings)
# will be written back below in a single pass so every row/column lines up
# exactly with the new row numbering.
$ws.Rows("2:17").ClearContents()

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $item.Title
    $ws.Cells.Item($r, 3).Value = $item.Category
    $ws.Cells.Item($r, 4).Value = $item.Price
    $ws.Cells.Item($r, 5).Value = $item.Deadline
    $ws.Cells.Item($r, 6).Value = $item.Url
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
    $ws.Cells.Item($r, 7).Value = $item.Score
    if ($item.Skill -ne "") {
        $ws.Cells.Item($r, 8).Value = $item.Skill
    } else {
        $ws.Cells.Item($r, 8).Value = ""
    }

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $item.Url) | Out-Null
}

# Widen column B (title) and column D (price) to fit the new, longer listing text.
$ws.Columns.Item(2).ColumnWidth = 54 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 32 - (5/6)

